$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H column (magnitude counts) gets a new number format (#,##0) + Arial 13 font,
#     matching the "totals" font already used elsewhere (fontId 1). Apply in two
#     contiguous blocks so the currently-empty H5 cell (which holds no H value,
#     I5 holds the value instead) is not touched / materialised.
$r1 = $ws.Range("H1:H4")
$r1.NumberFormat = "#,##0"
$r1.Font.Name = "Arial"
$r1.Font.Size = 13

$r2 = $ws.Range("H6:H25")
$r2.NumberFormat = "#,##0"
$r2.Font.Name = "Arial"
$r2.Font.Size = 13

# I5 (total records value) gets the same treatment.
$r3 = $ws.Range("I5")
$r3.NumberFormat = "#,##0"
$r3.Font.Name = "Arial"
$r3.Font.Size = 13

# Row 26 used to hold literal placeholder text ("=SUM(H1:H25)" / "=sum(H1:H25")
# in H26/I26. Replace them with real formulas, and add the new J26 ratio cell.
$ws.Range("H26").Formula = "=SUM(H1:H25)"
$ws.Range("H26").NumberFormat = "#,##0"

$ws.Range("I26").Formula = "=SUM(I1:I25)"
$ws.Range("I26").NumberFormat = "#,##0"

$ws.Range("J26").Formula = "=H26/I26"
$ws.Range("J26").NumberFormat = "0.00%"

# Match the author's final selection (cell J26, the new ratio cell).
[void]$ws.Range("J26").Select()
